$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.234.48'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.833.35'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.56'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6216'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07373'
$ws.Range('E8').Value = '  -1.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2906'
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.28'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07676'
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = '1.829.99'
$ws.Range('E12').Value = '  -0.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.978'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6699'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.76'
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008978'
$ws.Range('E16').Value = '  -3.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.870'
$ws.Range('E17').Value = '  -1.47%  '
$ws.Range('D18').Value = '29.204.97'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').Value = '2.068.85'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '236.75'
$ws.Range('E20').Value = '  +1.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.51'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.351'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9987'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.25'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1404'
$ws.Range('E26').Value = '  +0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.548'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.489'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05765'
$ws.Range('E30').Value = '  +3.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.107'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.091'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.209'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.873'
$ws.Range('E34').Value = '  +1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7326'
$ws.Range('E35').Value = '  -2.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.143'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.601'
$ws.Range('E37').Value = '  -2.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.862'
$ws.Range('E38').Value = '  +3.05%  '
$ws.Range('D39').Value = '1.224.77'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01754'
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.270'
$ws.Range('E41').Value = '  -3.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9068'
$ws.Range('E42').Value = '  +0.84%  '
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.60'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '1.974.48'
$ws.Range('E45').Value = '  -2.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.45'
$ws.Range('E46').Value = '  -0.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5039'
$ws.Range('E47').Value = '  -1.42%  '
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4028'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.132'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1136'
$ws.Range('E51').Value = '  +3.13%  '
